$wb = $excel.ActiveWorkbook

# --- Rename header labels on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$data = @(
    @(45060.99999999999, 20, 19.99999997272044, 20.00000002498298),
    @(45067.99999999999, 20, 19.99999997240842, 20.00000002521156),
    @(45074.99999999999, 20, 19.9999998469148, 20.0000001530731),
    @(45081.99999999999, 20, 19.99999948858085, 20.00000051682872),
    @(45088.99999999999, 20, 19.99999903082195, 20.00000098728876),
    @(45095.99999999999, 20, 19.99999846485823, 20.0000015356406),
    @(45102.99999999999, 20, 19.99999786538467, 20.0000022497866),
    @(45109.99999999999, 20, 19.99999722192781, 20.00000294745938),
    @(45116.99999999999, 20, 19.99999646866586, 20.00000371126722),
    @(45123.99999999999, 20, 19.99999565973303, 20.00000454750002)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# --- Match formatting to the sibling sheets: bold/centered header row and
#     the shared date format on column A ---
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)

$wsWeekly.Activate()
$wsWeekly.Range("A1").Select() | Out-Null
